$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Reprogramación" column (D): fill in the rescheduled class dates ---
# D14 currently holds the placeholder text "fecha actual"; replace it with
# the actual reprogrammed date (same day as A14) using the same date
# formatting/style as the neighbouring date cells in column A.
$ws.Range("A14").Copy()
$ws.Range("D14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D14").Value2 = 45155

$ws.Range("A15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value2 = 45157

$ws.Range("A16").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value2 = 45164

$ws.Range("A17").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value2 = 45169

$ws.Range("A18").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value2 = 45171

# --- Swap the two footer notes (B23 / B24) ---
$b23 = $ws.Range("B23").Value2
$b24 = $ws.Range("B24").Value2
$ws.Range("B23").Value2 = $b24
$ws.Range("B24").Value2 = $b23

# --- Update the active selection to D18 (last edited cell) ---
$ws.Range("D18").Select()
